$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for the greedy algorithm results (columns C: RoCoF, D: AFV, E: AV)
$data = @{
    2  = @{ C = 0.053;              D = 0.295;               E = 0.295 }
    3  = @{ C = 0.07066666666666667; D = 0.18;                E = 0.1833333333333333 }
    4  = @{ C = 0;                  D = 0.1776666666666667;  E = 0.1776666666666667 }
    5  = @{ C = 0.203;              D = 0.2553333333333334;  E = 0.3316666666666667 }
    6  = @{ C = 0.1963333333333333; D = 0.297;               E = 0.355 }
    7  = @{ C = 0.05233333333333334;D = 0.2276666666666667;  E = 0.2306666666666667 }
    8  = @{ C = 0.03633333333333334;D = 0.221;               E = 0.222 }
    9  = @{ C = 0.041;              D = 0.253;               E = 0.2533333333333334 }
    10 = @{ C = 0.001;              D = 0.2246666666666667;  E = 0.2246666666666667 }
    11 = @{ C = 0;                  D = 0.04666666666666667; E = 0.04666666666666667 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
}
